$wb = $excel.ActiveWorkbook

# zh-cn sheet: "Correspond Handoff Datetime" (E) and "Correspond Handback DateTime" (H)
# for the first two data rows move forward in time (report regenerated).
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E2").Value = "2016-03-11 20:14:48"
$wsZhCn.Range("E3").Value = "2016-03-11 20:14:48"
$wsZhCn.Range("H2").Value = "2016-03-11 20:15:09"
$wsZhCn.Range("H3").Value = "2016-03-11 20:15:09"

# de-de sheet: same columns, same two rows.
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E2").Value = "2016-03-11 20:14:51"
$wsDeDe.Range("E3").Value = "2016-03-11 20:14:51"
$wsDeDe.Range("H2").Value = "2016-03-11 20:15:15"
$wsDeDe.Range("H3").Value = "2016-03-11 20:15:15"
